# Weekly update: insert a new price record (row 296) for
# "Vega Modelo de Temuco - Ciboulette", shifting the existing
# historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 296 (shifts 296:317 -> 297:318,
# preserving all existing cell values/styles/formatting).
$ws.Rows.Item(296).Insert()

# Populate the newly inserted row 296 with the new weekly record.
$ws.Cells.Item(296, 1).Value = 10
$ws.Cells.Item(296, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(296, 3).Value = "La Araucanía"
$ws.Cells.Item(296, 4).Value = 44931
$ws.Cells.Item(296, 5).Value = 9
$ws.Cells.Item(296, 6).Value = 100112039
$ws.Cells.Item(296, 7).Value = "Ciboulette"
$ws.Cells.Item(296, 8).Value = "Sin especificar"
$ws.Cells.Item(296, 9).Value = "Primera"
$ws.Cells.Item(296, 10).Value = 80
$ws.Cells.Item(296, 11).Value = 5000
$ws.Cells.Item(296, 12).Value = 5000
$ws.Cells.Item(296, 13).Value = 5000
$ws.Cells.Item(296, 14).Value = "$/docena de atados"
$ws.Cells.Item(296, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(296, 16).Value = 1667
$ws.Cells.Item(296, 17).Value = 3
$ws.Cells.Item(296, 18).Value = "Hortaliza"
